$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend two leading spaces to the answer letters in column B, rows 2-21
for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "  " + $cell.Text
}

# Remove the now-unused extra question rows 22-41
$ws.Rows("22:41").Delete()
